$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 13.7015983672966
$ws.Range("C2").Value = 7.027492222474772
$ws.Range("D2").Value = 9.216955118608999
$ws.Range("E2").Value = 13.49660647338201
$ws.Range("F2").Value = 31.84037344995379
$ws.Range("J2").Value = 9.917056146039144
$ws.Range("M2").Value = 16.5053817125972
$ws.Range("N2").Value = 17.78526280807522
$ws.Range("O2").Value = 23.84670232787552

$ws.Range("B3").Value = 13.2392116101996
$ws.Range("C3").Value = 6.614749596745622
$ws.Range("D3").Value = 9.207031294696167
$ws.Range("E3").Value = 13.51183840921222
$ws.Range("F3").Value = 31.85815364066598
$ws.Range("J3").Value = 9.941423444315424
$ws.Range("M3").Value = 16.3710400245124
$ws.Range("N3").Value = 17.84211868669211
$ws.Range("O3").Value = 23.89158391864082

$ws.Range("B4").Value = 12.9489857401828
$ws.Range("C4").Value = 6.347681334806301
$ws.Range("D4").Value = 9.202121101417852
$ws.Range("E4").Value = 13.52329630282917
$ws.Range("F4").Value = 31.87758136147763
$ws.Range("J4").Value = 9.957505084942829
$ws.Range("M4").Value = 16.29064116373498
$ws.Range("N4").Value = 17.87886909178367
$ws.Range("O4").Value = 23.92519417823694

$ws.Range("B5").Value = 12.82931357957181
$ws.Range("C5").Value = 6.235484225123206
$ws.Range("D5").Value = 9.20041941785721
$ws.Range("E5").Value = 13.52849510792066
$ws.Range("F5").Value = 31.88763557887668
$ws.Range("J5").Value = 9.964340431528107
$ws.Range("M5").Value = 16.25843208630279
$ws.Range("N5").Value = 17.89430915444411
$ws.Range("O5").Value = 23.94040852223315

$ws.Range("B6").Value = 12.80936313616823
$ws.Range("C6").Value = 6.21665299767611
$ws.Range("D6").Value = 9.200154978009124
$ws.Range("E6").Value = 13.52939035696503
$ws.Range("F6").Value = 31.88943405396411
$ws.Range("J6").Value = 9.965492473678365
$ws.Range("M6").Value = 16.2531180949978
$ws.Range("N6").Value = 17.89690102127192
$ws.Range("O6").Value = 23.9430263989741

$ws.Range("B7").Value = 12.94737722169324
$ws.Range("C7").Value = 6.346181730231794
$ws.Range("D7").Value = 9.202096937991827
$ws.Range("E7").Value = 13.52336427113694
$ws.Range("F7").Value = 31.87770830710362
$ws.Range("J7").Value = 9.957596126811834
$ws.Range("M7").Value = 16.29020449942117
$ws.Range("N7").Value = 17.87907544176335
$ws.Range("O7").Value = 23.9253932241733

$ws.Range("B8").Value = 13.54358551243872
$ws.Range("C8").Value = 6.888045114935079
$ws.Range("D8").Value = 9.213288936087732
$ws.Range("E8").Value = 13.50142160930455
$ws.Range("F8").Value = 31.84473635418416
$ws.Range("J8").Value = 9.925225716085166
$ws.Range("M8").Value = 16.45864523580779
$ws.Range("N8").Value = 17.80448536051644
$ws.Range("O8").Value = 23.86091890497133

$ws.Range("B9").Value = 14.65525163213465
$ws.Range("C9").Value = 7.879333994089476
$ws.Range("D9").Value = 9.244541137733558
$ws.Range("E9").Value = 13.47508629247726
$ws.Range("F9").Value = 31.84768058907767
$ws.Range("J9").Value = 9.870620505705677
$ws.Range("M9").Value = 16.80419742723289
$ws.Range("N9").Value = 17.67276802438999
$ws.Range("O9").Value = 23.78267576154558

$ws.Range("B10").Value = 15.42820276905078
$ws.Range("C10").Value = 8.539740062976282
$ws.Range("D10").Value = 9.273056944259714
$ws.Range("E10").Value = 13.46589525205472
$ws.Range("F10").Value = 31.89108226027635
$ws.Range("J10").Value = 9.835892442084591
$ws.Range("M10").Value = 17.06554746032982
$ws.Range("N10").Value = 17.58479582089568
$ws.Range("O10").Value = 23.75476397642571

$ws.Range("B11").Value = 15.76878892440213
$ws.Range("C11").Value = 8.822712928729572
$ws.Range("D11").Value = 9.287207774848188
$ws.Range("E11").Value = 13.46391317374796
$ws.Range("F11").Value = 31.91975981593065
$ws.Range("J11").Value = 9.821260273851715
$ws.Range("M11").Value = 17.18566852561821
$ws.Range("N11").Value = 17.54667092570324
$ws.Range("O11").Value = 23.74851803714995

$ws.Range("B12").Value = 15.89606419454962
$ws.Range("C12").Value = 8.927365604990435
$ws.Range("D12").Value = 9.292733219767321
$ws.Range("E12").Value = 13.46347804219943
$ws.Range("F12").Value = 31.93189976828977
$ws.Range("J12").Value = 9.815886801823517
$ws.Range("M12").Value = 17.23129868414038
$ws.Range("N12").Value = 17.53250530245582
$ws.Range("O12").Value = 23.74708198612878

$ws.Range("B13").Value = 15.86873030120427
$ws.Range("C13").Value = 8.904937991620363
$ws.Range("D13").Value = 9.291535840641419
$ws.Range("E13").Value = 13.46355773998241
$ws.Range("F13").Value = 31.92922835016694
$ws.Range("J13").Value = 9.817036633058436
$ws.Range("M13").Value = 17.22146560087797
$ws.Range("N13").Value = 17.53554406464113
$ws.Range("O13").Value = 23.74734992456337

$ws.Range("B14").Value = 15.77929454209845
$ws.Range("C14").Value = 8.831372892446671
$ws.Range("D14").Value = 9.287659028823937
$ws.Range("E14").Value = 13.46387105822365
$ws.Range("F14").Value = 31.92073294588634
$ws.Range("J14").Value = 9.820814841175928
$ws.Range("M14").Value = 17.18941987769655
$ws.Range("N14").Value = 17.54550007635286
$ws.Range("O14").Value = 23.74838126465517

$ws.Range("B15").Value = 15.72428845358634
$ws.Range("C15").Value = 8.785986380176462
$ws.Range("D15").Value = 9.285306014120403
$ws.Range("E15").Value = 13.46410402904832
$ws.Range("F15").Value = 31.9156958579211
$ws.Range("J15").Value = 9.823150898063702
$ws.Range("M15").Value = 17.16980854460867
$ws.Range("N15").Value = 17.55163374443551
$ws.Range("O15").Value = 23.74913402581262

$ws.Range("B16").Value = 15.40571287667431
$ws.Range("C16").Value = 8.520896501402243
$ws.Range("D16").Value = 9.272155649351999
$ws.Range("E16").Value = 13.46606897145912
$ws.Range("F16").Value = 31.88938755180819
$ws.Range("J16").Value = 9.836872128018671
$ws.Range("M16").Value = 17.05771912746834
$ws.Range("N16").Value = 17.58732540726868
$ws.Range("O16").Value = 23.75530211896327

$ws.Range("B17").Value = 15.20736817297324
$ws.Range("C17").Value = 8.353807336579004
$ws.Range("D17").Value = 9.264388392213791
$ws.Range("E17").Value = 13.46783710108334
$ws.Range("F17").Value = 31.87553382904624
$ws.Range("J17").Value = 9.845588084620179
$ws.Range("M17").Value = 16.98924703479068
$ws.Range("N17").Value = 17.60970554788352
$ws.Range("O17").Value = 23.76073961165904

$ws.Range("B18").Value = 15.09225279293164
$ws.Range("C18").Value = 8.256059224662994
$ws.Range("D18").Value = 9.260031999635391
$ws.Range("E18").Value = 13.46906111426412
$ws.Range("F18").Value = 31.86840678231346
$ws.Range("J18").Value = 9.850711016601755
$ws.Range("M18").Value = 16.94998167635885
$ws.Range("N18").Value = 17.62275634831692
$ws.Range("O18").Value = 23.76447428955935

$ws.Range("B19").Value = 15.05310291956915
$ws.Range("C19").Value = 8.222681063392177
$ws.Range("D19").Value = 9.258576165841642
$ws.Range("E19").Value = 13.469511123422
$ws.Range("F19").Value = 31.86613828030121
$ws.Range("J19").Value = 9.85246440904166
$ws.Range("M19").Value = 16.93670840486606
$ws.Range("N19").Value = 17.62720578287468
$ws.Range("O19").Value = 23.76584301547683

$ws.Range("B20").Value = 15.22859003542071
$ws.Range("C20").Value = 8.371764250167772
$ws.Range("D20").Value = 9.265203747641506
$ws.Range("E20").Value = 13.46762745883912
$ws.Range("F20").Value = 31.87692154217524
$ws.Range("J20").Value = 9.844648899701797
$ws.Range("M20").Value = 16.9965240432578
$ws.Range("N20").Value = 17.6073046933951
$ws.Range("O20").Value = 23.76009792820835

$ws.Range("B21").Value = 15.80561087495616
$ws.Range("C21").Value = 8.853048622011949
$ws.Range("D21").Value = 9.288793236726073
$ws.Range("E21").Value = 13.46377047492545
$ws.Range("F21").Value = 31.92319354395919
$ws.Range("J21").Value = 9.819700548144985
$ws.Range("M21").Value = 17.19882887701678
$ws.Range("N21").Value = 17.54256839461802
$ws.Range("O21").Value = 23.74805310963621

$ws.Range("B22").Value = 16.17279174590522
$ws.Range("C22").Value = 9.153009879229762
$ws.Range("D22").Value = 9.305181416562613
$ws.Range("E22").Value = 13.46308788232171
$ws.Range("F22").Value = 31.96089505694698
$ws.Range("J22").Value = 9.804371062511166
$ws.Range("M22").Value = 17.33186400314047
$ws.Range("N22").Value = 17.50184127056324
$ws.Range("O22").Value = 23.7455971570205

$ws.Range("B23").Value = 15.97776276848828
$ws.Range("C23").Value = 8.994247231359996
$ws.Range("D23").Value = 9.296346804696052
$ws.Range("E23").Value = 13.46328428979962
$ws.Range("F23").Value = 31.94009223164537
$ws.Range("J23").Value = 9.812463495210787
$ws.Range("M23").Value = 17.26079715953222
$ws.Range("N23").Value = 17.52343366830405
$ws.Range("O23").Value = 23.74641205732102

$ws.Range("B24").Value = 15.21899900040454
$ws.Range("C24").Value = 8.3636511825735
$ws.Range("D24").Value = 9.26483478540492
$ws.Range("E24").Value = 13.46772159172906
$ws.Range("F24").Value = 31.87629154752479
$ws.Range("J24").Value = 9.845073156273454
$ws.Range("M24").Value = 16.99323379312553
$ws.Range("N24").Value = 17.6083895459739
$ws.Range("O24").Value = 23.76038613765745

$ws.Range("B25").Value = 14.36163465505799
$ws.Range("C25").Value = 7.620626358039075
$ws.Range("D25").Value = 9.235102014581646
$ws.Range("E25").Value = 13.48042503715313
$ws.Range("F25").Value = 31.83963885888553
$ws.Range("J25").Value = 9.884444702569306
$ws.Range("M25").Value = 16.70927034934629
$ws.Range("N25").Value = 17.70685048885852
$ws.Range("O25").Value = 23.74734992456337
